# "Generate Report for handoff"
#
# The localization-status report is regenerated: the failed-handoff row
# (367ec86e.../"Handoff transform failed") is gone, replaced by a fresh
# successful handoff for a new source file (a8ab76e9-...), with new
# handoff package file names and new handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid  = "367ec86e-8185-4d85-b08e-47192c1deefb"
$newGuid  = "a8ab76e9-1fdd-4c9a-a605-f33bbd8e11cb"
$oldHash  = "5ec611ce3611e9cd550f0af0fc14ef229e79e7fb"
$newHash  = "5c67c3533df65061b6d2c008b9de67fb3c3ac9a1"

$hlFontColor = 15570276   # RGB(0x64,0x95,0xED) == FF6495ED, the workbook's HyperLink font color

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/da4dd7ef2f6aeeb795657c54f341d91f8918d85e"

# ---- Sheet 1: Overview -----------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Row 3 was the failed-handoff row (921e6f04.../"Handoff transform failed");
# it is removed and the config row moves up to become row 3.
$ws.Rows(3).Delete()

$ws.Cells.Replace($oldGuid, $newGuid)

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + "/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + "/.localization-config", "", "", ".localization-config")

$hr = $ws.Range("A2:A3")
$hr.Font.Underline = 2
$hr.Font.Color = $hlFontColor

# ---- Sheet 2: zh-cn -----------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Rows(3).Delete()

$ws.Cells.Replace($oldGuid, $newGuid)
$ws.Cells.Replace($oldHash, $newHash)
$ws.Range("D2").Value = "2016-01-08 20:15:20"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + "/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3617260ce96b8f558a6b4b84a88b1b47e479a91/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/" + $newGuid + "." + $newHash + ".zh-cn.xlf", "", "", $newGuid + "." + $newHash + ".zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + "/.localization-config", "", "", ".localization-config")

$hr = $ws.Range("A2:A3,C2")
$hr.Font.Underline = 2
$hr.Font.Color = $hlFontColor

# ---- Sheet 3: de-de -----------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Rows(3).Delete()

$ws.Cells.Replace($oldGuid, $newGuid)
$ws.Cells.Replace($oldHash, $newHash)
$ws.Range("D2").Value = "2016-01-08 20:15:28"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + "/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e3d38df4b59fea4bfb2c38f7793b3f9779280b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/" + $newGuid + "." + $newHash + ".de-de.xlf", "", "", $newGuid + "." + $newHash + ".de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + "/.localization-config", "", "", ".localization-config")

$hr = $ws.Range("A2:A3,C2")
$hr.Font.Underline = 2
$hr.Font.Color = $hlFontColor
